# Apply updated cryptocurrency price/volume data to the sheet.
# Force-text (leading apostrophe) is used for Price (column D) values
# that would otherwise be auto-parsed as numbers by Excel, so they keep
# the exact original text representation (e.g. trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.853.61"
$ws.Range("E2").Value = "  -4.00%  "
$ws.Range("D3").Value = "2.965.78"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'557.24"
$ws.Range("E5").Value = "  -3.80%  "
$ws.Range("D6").Value = "'133.63"
$ws.Range("E6").Value = "  +6.16%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +3.39%  "
$ws.Range("D9").Value = "2.959.38"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  -4.89%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'32.99"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "3.452.90"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  +9.47%  "
$ws.Range("D18").Value = "2.963.67"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "57.833.68"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("D20").Value = "'418.90"
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "'6.98"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "'13.02"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "'79.58"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "'7.58"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("E30").Value = "  +6.22%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'6.07"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "'0.1000"
$ws.Range("E33").Value = "  +7.15%  "
$ws.Range("D34").Value = "'5.66"
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "'2.13"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "'0.938"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "0.0₃0691"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").Value = "'48.53"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "'8.54"
$ws.Range("E39").Value = "  +6.53%  "
$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0351"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'381.53"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "2.681.69"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "'122.38"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").Value = "'2.00"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "'23.54"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -0.18%  "
